$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Abril de 2020 a las 02:22"

# Row 4: Estados Unidos - refreshed totals
$ws.Range("B4").Value = 960525
$ws.Range("C4").Value = 35293
$ws.Range("D4").Value = 118162
$ws.Range("E4").Value = 788115
$ws.Range("F4").Value = 15110
$ws.Range("G4").Value = 2055
$ws.Range("H4").Value = 54248

# Row 14: Brasil - refreshed totals
$ws.Range("B14").Value = 59196
$ws.Range("C14").Value = 6201
$ws.Range("D14").Value = 29160
$ws.Range("E14").Value = 25991
$ws.Range("F14").Value = 8318
$ws.Range("G14").Value = 375
$ws.Range("H14").Value = 4045

# Japon's case counts jumped enough to overtake Mexico and Chile in the
# sorted ranking, so it now sits right after Austria (row 28), pushing
# Mexico to row 29 and Chile to row 30. Their own totals are unchanged.
$ws.Range("A28").Value = "Japon"
$ws.Range("B28").Value = 13231
$ws.Range("C28").Value = 519
$ws.Range("D28").Value = 1656
$ws.Range("E28").Value = 11215
$ws.Range("F28").Value = 287
$ws.Range("G28").Value = 15
$ws.Range("H28").Value = 360

$ws.Range("A29").Value = "Mexico"
$ws.Range("B29").Value = 12872
$ws.Range("C29").Value = 1239
$ws.Range("D29").Value = 7149
$ws.Range("E29").Value = 4502
$ws.Range("F29").Value = 378
$ws.Range("G29").Value = 152
$ws.Range("H29").Value = 1221

$ws.Range("A30").Value = "Chile"
$ws.Range("B30").Value = 12858
$ws.Range("C30").Value = 552
$ws.Range("D30").Value = 6746
$ws.Range("E30").Value = 5931
$ws.Range("F30").Value = 418
$ws.Range("G30").Value = 7
$ws.Range("H30").Value = 181
